# Refresh the "articles" sheet with a new day's scraped news data and add
# a "썸네일" (thumbnail) column (D) with a thumbnail image URL per article.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "제목"
$ws.Range("B1").Value = "링크"
$ws.Range("C1").Value = "신문사"
$ws.Range("D1").Value = "썸네일"

$ws.Range("A2").Value = "“추석 이동 줄었더니…가정폭력·교통사고 줄었다”(종합)"
$ws.Range("B2").Value = "https://www.seoul.co.kr/news/newsView.php?id=20201004500093&wlog_tag3=naver"
$ws.Range("C2").Value = "네이버뉴스"
$ws.Range("D2").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F081%2F2020%2F10%2F04%2F3128560.jpg&type=ofullfill80_80_q75_re2"

$ws.Range("A3").Value = "`"이 시국에 장·차관 홍보하나`"...복지부 추석 포스터 논란"
$ws.Range("B3").Value = "https://biz.chosun.com/site/data/html_dir/2020/10/04/2020100400542.html?utm_source=naver&utm_medium=original&utm_campaign=biz"
$ws.Range("C3").Value = "네이버뉴스"
$ws.Range("D3").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F366%2F2020%2F10%2F04%2F598214.jpg&type=ofullfill80_80_q75_re2"

$ws.Range("A4").Value = "추석선물 백화점서 역대 가장 많이 샀다…`"비쌀수록 인기`"(종합)"
$ws.Range("B4").Value = "http://yna.kr/AKR20201004029351030?did=1195m"
$ws.Range("C4").Value = "네이버뉴스"
$ws.Range("D4").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F001%2F2020%2F10%2F04%2F11919710.jpg&type=ofullfill80_80_q75_re2"

$ws.Range("A5").Value = "복지부, '장관 얼굴' 추석 포스터 논란에 `"물의 일으켜 송구`""
$ws.Range("B5").Value = "http://www.busan.com/view/busan/view.php?code=2020100417233211300"
$ws.Range("C5").Value = "네이버뉴스"
$ws.Range("D5").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F082%2F2020%2F10%2F04%2F1033079.jpg&type=ofullfill80_80_q75_re2"

$ws.Range("A6").Value = "추석 연휴 이후에도 수도권 고위험시설 영업중단은 계속"
$ws.Range("B6").Value = "http://news.kmib.co.kr/article/view.asp?arcid=0015072804&code=61121111&cp=nv"
$ws.Range("C6").Value = "네이버뉴스"
$ws.Range("D6").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F005%2F2020%2F10%2F04%2F1367473.jpg&type=ofullfill80_80_q75_re2"

$ws.Range("A7").Value = "수월한 추석 귀경길…`"고속도로, 평소 주말보다 원활`""
$ws.Range("B7").Value = "http://www.newsis.com/view/?id=NISX20201004_0001185650&cID=10201&pID=10200"
$ws.Range("C7").Value = "네이버뉴스"
$ws.Range("D7").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F003%2F2020%2F10%2F04%2F10106636.jpg&type=ofullfill80_80_q75_re2"

$ws.Range("A8").Value = "추석 나훈아 소신발언에…野 `"속 시원하게 文 비판`" 與 `"오독 말라`""
$ws.Range("B8").Value = "http://news.mk.co.kr/newsRead.php?no=1012678&year=2020"
$ws.Range("C8").Value = "네이버뉴스"
$ws.Range("D8").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F009%2F2020%2F10%2F04%2F4667218.jpg&type=ofullfill80_80_q75_re2"

$ws.Range("A9").Value = "추석 전 못 받은 2차 재난지원금, 연휴 뒤 10~11월에 지급"
$ws.Range("B9").Value = "https://www.chosun.com/economy/2020/10/04/YTRRE6TK2FGEHI6TF2PDWEUGVE/?utm_source=naver&utm_medium=original&utm_campaign=news"
$ws.Range("C9").Value = "네이버뉴스"
$ws.Range("D9").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F023%2F2020%2F10%2F04%2F3566076.jpg&type=ofullfill80_80_q75_re2"

$ws.Range("A10").Value = "코로나19 신규확진 64명 추석 귀성·귀경객 2명 확진 판정"
$ws.Range("B10").Value = "http://www.dongascience.com/news/view/40274"
$ws.Range("C10").Value = "네이버뉴스"
$ws.Range("D10").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F584%2F2020%2F10%2F04%2F10730.jpg&type=ofullfill80_80_q75_re2"

# Note: the title starts with a literal apostrophe. Assigning a string
# that begins with "'" makes Excel treat it as a text quote-prefix (it
# strips the leading apostrophe from the stored value and flags the cell
# as quote-prefixed). Doubling the apostrophe keeps one literal leading
# apostrophe in the stored value, and ClearFormats() afterwards drops the
# quote-prefix flag again, leaving a plain cell with the exact text.
$ws.Range("A11").Value = "''코로나19' 덮친 추석 연휴, 교통사고·112신고 줄었다"
$ws.Range("A11").ClearFormats()
$ws.Range("B11").Value = "https://www.nocutnews.co.kr/news/5422084"
$ws.Range("C11").Value = "네이버뉴스"
$ws.Range("D11").Value = "https://search.pstatic.net/common/?src=https%3A%2F%2Fimgnews.pstatic.net%2Fimage%2Forigin%2F079%2F2020%2F10%2F04%2F3413806.jpg&type=ofullfill80_80_q75_re2"
